# Added capex functionality (biomass -> gas -> km chain via new hub nodes).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sources sheet: add a new Biomass (BM) source row
# ---------------------------------------------------------------------------
$wsSources = $wb.Worksheets.Item("Sources")
$wsSources.Cells.Item(4, 1).Value = "BM"
$wsSources.Cells.Item(4, 2).Value = 0
$wsSources.Cells.Item(4, 3).Value = 0.02
$wsSources.Cells.Item(4, 4).Value = "biomass"
$wsSources.Cells.Item(4, 5).Value = 0.04

# ---------------------------------------------------------------------------
# Sinks sheet: rename the Gasoline demand sink to a Kilometers (km) sink
# ---------------------------------------------------------------------------
$wsSinks = $wb.Worksheets.Item("Sinks")
$wsSinks.Cells.Item(2, 1).Value = "Kilometers"
$wsSinks.Cells.Item(2, 4).Value = "km"

# ---------------------------------------------------------------------------
# Transformers sheet: Refinery now has a capex, plus two new transformers
# (Gtkm: gasoline -> km, B2gas: biomass -> gasoline)
# ---------------------------------------------------------------------------
$wsTransformers = $wb.Worksheets.Item("Transformers")
$wsTransformers.Cells.Item(2, 3).Value = 200

$wsTransformers.Cells.Item(4, 1).Value = "Gtkm"
$wsTransformers.Cells.Item(4, 2).Value = "gasoline"
$wsTransformers.Cells.Item(4, 3).Value = 0
$wsTransformers.Cells.Item(4, 4).Value = 0
$wsTransformers.Cells.Item(4, 5).Value = 0.4
$wsTransformers.Cells.Item(4, 6).Value = "km"
$wsTransformers.Cells.Item(4, 7).Value = 1

$wsTransformers.Cells.Item(5, 1).Value = "B2gas"
$wsTransformers.Cells.Item(5, 2).Value = "biomass"
$wsTransformers.Cells.Item(5, 3).Value = 0
$wsTransformers.Cells.Item(5, 4).Value = 0
$wsTransformers.Cells.Item(5, 5).Value = 0.5
$wsTransformers.Cells.Item(5, 6).Value = "gasoline"
$wsTransformers.Cells.Item(5, 7).Value = 1

# ---------------------------------------------------------------------------
# Connectors sheet: route old ref2gas/mtg2gas connectors into a new GasHub,
# then add the new GasHub -> Gtkm -> KmHub -> Kilometers chain plus the
# biomass -> B2gas -> GasHub chain
# ---------------------------------------------------------------------------
$wsConnectors = $wb.Worksheets.Item("Connectors")
$wsConnectors.Cells.Item(4, 3).Value = "GasHub"
$wsConnectors.Cells.Item(5, 3).Value = "GasHub"

$wsConnectors.Cells.Item(6, 1).Value = "gas2km"
$wsConnectors.Cells.Item(6, 2).Value = "GasHub"
$wsConnectors.Cells.Item(6, 3).Value = "Gtkm"
$wsConnectors.Cells.Item(6, 4).Value = "gasoline"

$wsConnectors.Cells.Item(7, 1).Value = "kmtohub"
$wsConnectors.Cells.Item(7, 2).Value = "Gtkm"
$wsConnectors.Cells.Item(7, 3).Value = "KmHub"
$wsConnectors.Cells.Item(7, 4).Value = "km"

$wsConnectors.Cells.Item(8, 1).Value = "hub2sink"
$wsConnectors.Cells.Item(8, 2).Value = "KmHub"
$wsConnectors.Cells.Item(8, 3).Value = "Kilometers"
$wsConnectors.Cells.Item(8, 4).Value = "km"

$wsConnectors.Cells.Item(9, 1).Value = "bm2btg"
$wsConnectors.Cells.Item(9, 2).Value = "BM"
$wsConnectors.Cells.Item(9, 3).Value = "B2gas"
$wsConnectors.Cells.Item(9, 4).Value = "biomass"

$wsConnectors.Cells.Item(10, 1).Value = "btg2gas"
$wsConnectors.Cells.Item(10, 2).Value = "B2gas"
$wsConnectors.Cells.Item(10, 3).Value = "GasHub"
$wsConnectors.Cells.Item(10, 4).Value = "gasoline"

# ---------------------------------------------------------------------------
# Hubs sheet: header rename + two new hub rows (GasHub, KmHub)
# ---------------------------------------------------------------------------
$wsHubs = $wb.Worksheets.Item("Hubs")
$wsHubs.Cells.Item(1, 2).Value = "EnergyType"

$wsHubs.Cells.Item(2, 1).Value = "GasHub"
$wsHubs.Cells.Item(2, 2).Value = "gasoline"
$wsHubs.Cells.Item(2, 3).Value = 0
$wsHubs.Cells.Item(2, 4).Value = 0

$wsHubs.Cells.Item(3, 1).Value = "KmHub"
$wsHubs.Cells.Item(3, 2).Value = "km"
$wsHubs.Cells.Item(3, 3).Value = 0
$wsHubs.Cells.Item(3, 4).Value = 0

# ---------------------------------------------------------------------------
# View state: per-sheet selection + which sheet/cell is active overall.
# Transformers ends up the active sheet (activate it last).
# ---------------------------------------------------------------------------
$wsSources.Activate()
$wsSources.Range("C39").Select()

$wsSinks.Activate()
$wsSinks.Range("A39").Select()

$wsConnectors.Activate()
$wsConnectors.Range("D10").Select()

$wsHubs.Activate()
$wsHubs.Range("E3").Select()

$wsTransformers.Activate()
$wsTransformers.Range("D9").Select()
